$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 21
$ws.Range("H21").Value = 10000
$ws.Range("J21").Value = 10000
$ws.Range("L21").Value = 10000
$ws.Range("N21").Value = -10936
# Row 23
$ws.Range("H23").Value = 10000
$ws.Range("J23").Value = 10000
$ws.Range("L23").Value = 10000
$ws.Range("N23").Value = -10468
# Row 41
$ws.Range("H41").Value = 334.16666
$ws.Range("I41").Value = 386.7
$ws.Range("K41").Value = 386.7
$ws.Range("M41").Value = 53.30000000000001
# Row 51
$ws.Range("H51").Value = 11345.25
$ws.Range("I51").Value = 11398
$ws.Range("J51").Value = 11292.5
$ws.Range("K51").Value = 11398
$ws.Range("L51").Value = 11292.5
$ws.Range("M51").Value = -10914
$ws.Range("N51").Value = -12260.5
# Row 57
$ws.Range("H57").Value = 90626.5
$ws.Range("J57").Value = 90626.5
$ws.Range("L57").Value = 271879.5
$ws.Range("N57").Value = -272877.5
# Row 86
$ws.Range("H86").Value = 4064.2307
$ws.Range("I86").Value = 3538.4443
$ws.Range("K86").Value = 3538.4443
$ws.Range("M86").Value = -2415.4443
# Row 89
$ws.Range("H89").Value = 4064.2307
$ws.Range("I89").Value = 3538.4443
$ws.Range("K89").Value = 17692.2215
$ws.Range("M89").Value = -12076.2215
# Row 92
$ws.Range("H92").Value = 717.6
$ws.Range("I92").Value = 582.13336
$ws.Range("J92").Value = 920.8
$ws.Range("K92").Value = 582.13336
$ws.Range("L92").Value = 920.8
$ws.Range("M92").Value = 665.86664
$ws.Range("N92").Value = -3416.8

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 39
$ws.Range("H39").Value = 1500
$ws.Range("I39").Value = 1500
$ws.Range("K39").Value = 1500
$ws.Range("M39").Value = -980
# Row 44
$ws.Range("H44").Value = 49999
$ws.Range("J44").Value = 49999
$ws.Range("L44").Value = 49999
$ws.Range("N44").Value = -50975
# Row 45
$ws.Range("H45").Value = 3432.6667
$ws.Range("I45").Value = 3106.6667
$ws.Range("K45").Value = 3106.6667
$ws.Range("M45").Value = -2729.6667
# Row 55
$ws.Range("H55").Value = 43723.75
$ws.Range("J55").Value = 49999
$ws.Range("L55").Value = 49999
$ws.Range("N55").Value = -50629
# Row 102
$ws.Range("H102").Value = 3738.6
$ws.Range("I102").Value = 3819.5862
$ws.Range("K102").Value = 3819.5862
$ws.Range("M102").Value = -2197.5862
# Row 110
$ws.Range("H110").Value = 985
$ws.Range("I110").Value = 985
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 985
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 1060
$ws.Range("N110").ClearContents()

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 82
$ws.Range("H82").Value = 42561
$ws.Range("J82").Value = 89977.664
$ws.Range("L82").Value = 89977.664
$ws.Range("N82").Value = -90743.664
# Row 85
$ws.Range("H85").Value = 42561
$ws.Range("J85").Value = 89977.664
$ws.Range("L85").Value = 89977.664
$ws.Range("N85").Value = -92629.664
# Row 94
$ws.Range("H94").Value = 4437.7104
$ws.Range("I94").Value = 1398.6923
$ws.Range("K94").Value = 1398.6923
$ws.Range("M94").Value = -947.6922999999999
# Row 132
$ws.Range("H132").Value = 74749.25
$ws.Range("J132").Value = 74749.25
$ws.Range("L132").Value = 74749.25
$ws.Range("N132").Value = -84869.25
# Row 133
$ws.Range("H133").Value = 93990
$ws.Range("J133").Value = 93990
$ws.Range("L133").Value = 93990
$ws.Range("N133").Value = -104110
# Row 134
$ws.Range("H134").Value = 8707.454
$ws.Range("I134").Value = 7907.5483
$ws.Range("K134").Value = 23722.6449
$ws.Range("M134").Value = -21187.6449

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 62
$ws.Range("H62").Value = 5881.3335
$ws.Range("I62").Value = 4000
$ws.Range("J62").Value = 7762.6665
$ws.Range("K62").Value = 4000
$ws.Range("L62").Value = 7762.6665
$ws.Range("M62").Value = -3376
$ws.Range("N62").Value = -9010.666499999999
# Row 65
$ws.Range("H65").Value = 5881.3335
$ws.Range("I65").Value = 4000
$ws.Range("J65").Value = 7762.6665
$ws.Range("K65").Value = 20000
$ws.Range("L65").Value = 38813.3325
$ws.Range("M65").Value = -16880
$ws.Range("N65").Value = -45053.3325
# Row 105
$ws.Range("H105").Value = 1618.1
$ws.Range("I105").Value = 1235.1666
$ws.Range("K105").Value = 1235.1666
$ws.Range("M105").Value = 511.8334
# Row 106
$ws.Range("H106").Value = 56988.5
$ws.Range("J106").Value = 56988.5
$ws.Range("L106").Value = 56988.5
$ws.Range("N106").Value = -59512.5
# Row 107
$ws.Range("H107").Value = 3159.35
$ws.Range("I107").Value = 3852.2666
$ws.Range("J107").Value = 1080.6
$ws.Range("K107").Value = 3852.2666
$ws.Range("L107").Value = 1080.6
$ws.Range("M107").Value = -1932.2666
$ws.Range("N107").Value = -4920.6
# Row 134
$ws.Range("H134").Value = 1845.4783
$ws.Range("I134").Value = 1910.0454
$ws.Range("J134").Value = 425
$ws.Range("K134").Value = 5730.1362
$ws.Range("L134").Value = 1275
$ws.Range("M134").Value = -3195.1362
$ws.Range("N134").Value = -6345

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 1906.6875
$ws.Range("J68").Value = 1906.6875
$ws.Range("L68").Value = 5720.0625
$ws.Range("N68").Value = -7342.0625
# Row 71
$ws.Range("H71").Value = 1906.6875
$ws.Range("J71").Value = 1906.6875
$ws.Range("L71").Value = 17160.1875
$ws.Range("N71").Value = -25272.1875
# Row 107
$ws.Range("H107").Value = 1223.875
$ws.Range("J107").Value = 1328.4166
$ws.Range("L107").Value = 3985.2498
$ws.Range("N107").Value = -7825.2498
# Row 132
$ws.Range("H132").Value = 3028.2173
$ws.Range("I132").Value = 1997
$ws.Range("J132").Value = 3691.1428
$ws.Range("K132").Value = 17973
$ws.Range("L132").Value = 33220.2852
$ws.Range("M132").Value = -15443
$ws.Range("N132").Value = -38280.2852
# Row 134
$ws.Range("H134").Value = 12837.944
$ws.Range("I134").Value = 3599.75
$ws.Range("J134").Value = 15477.429
$ws.Range("K134").Value = 10799.25
$ws.Range("L134").Value = 46432.287
$ws.Range("M134").Value = -5729.25
$ws.Range("N134").Value = -56572.287
# Row 139
$ws.Range("H139").Value = 11372.883
$ws.Range("I139").Value = 1017.5
$ws.Range("J139").Value = 20577.666
$ws.Range("K139").Value = 3052.5
$ws.Range("L139").Value = 61732.99800000001
$ws.Range("M139").Value = 2087.5
$ws.Range("N139").Value = -72012.99800000001

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 5696.846
$ws.Range("I70").Value = 5131.75
$ws.Range("J70").Value = 6601
$ws.Range("K70").Value = 5131.75
$ws.Range("L70").Value = 6601
$ws.Range("M70").Value = -4861.75
$ws.Range("N70").Value = -7141
# Row 73
$ws.Range("H73").Value = 5696.846
$ws.Range("I73").Value = 5131.75
$ws.Range("J73").Value = 6601
$ws.Range("K73").Value = 5131.75
$ws.Range("L73").Value = 6601
$ws.Range("M73").Value = -4195.75
$ws.Range("N73").Value = -8473
# Row 97
$ws.Range("H97").Value = 3166.8572
$ws.Range("I97").Value = 2243.25
$ws.Range("J97").Value = 4398.3335
$ws.Range("K97").Value = 2243.25
$ws.Range("L97").Value = 4398.3335
$ws.Range("M97").Value = -1747.25
$ws.Range("N97").Value = -5390.3335
# Row 135
$ws.Range("H135").Value = 82498.664
$ws.Range("J135").Value = 82498.664
$ws.Range("L135").Value = 82498.664
$ws.Range("N135").Value = -92638.664
# Row 138
$ws.Range("H138").Value = 97999
$ws.Range("J138").Value = 97999
$ws.Range("L138").Value = 97999
$ws.Range("N138").Value = -108279

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 1730
$ws.Range("I61").Value = 1778.0968
$ws.Range("K61").Value = 1778.0968
$ws.Range("M61").Value = -1576.0968
# Row 68
$ws.Range("H68").Value = 2510.8333
$ws.Range("I68").Value = 2510.8333
$ws.Range("K68").Value = 2510.8333
$ws.Range("M68").Value = -1761.8333
# Row 71
$ws.Range("H71").Value = 2510.8333
$ws.Range("I71").Value = 2510.8333
$ws.Range("K71").Value = 12554.1665
$ws.Range("M71").Value = -8810.166499999999
# Row 82
$ws.Range("H82").Value = 3287.7036
$ws.Range("I82").Value = 2029.8235
$ws.Range("K82").Value = 2029.8235
$ws.Range("M82").Value = -1668.8235
# Row 85
$ws.Range("H85").Value = 3287.7036
$ws.Range("I85").Value = 2029.8235
$ws.Range("K85").Value = 2029.8235
$ws.Range("M85").Value = -781.8235
# Row 87
$ws.Range("H87").Value = 132000
$ws.Range("J87").Value = 132000
$ws.Range("L87").Value = 132000
$ws.Range("N87").Value = -134246
# Row 90
$ws.Range("H90").Value = 132000
$ws.Range("J90").Value = 132000
$ws.Range("L90").Value = 396000
$ws.Range("N90").Value = -407232
# Row 100
$ws.Range("H100").Value = 1861.3334
$ws.Range("I100").Value = 1861.3334
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1861.3334
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -1320.3334
$ws.Range("N100").ClearContents()
# Row 113
$ws.Range("H113").Value = 1730
$ws.Range("I113").Value = 1778.0968
$ws.Range("K113").Value = 1778.0968
$ws.Range("M113").Value = 391.9032

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 3475.7273
$ws.Range("I122").Value = 3559.5557
$ws.Range("K122").Value = 10678.6671
$ws.Range("M122").Value = -8228.667099999999
# Row 132
$ws.Range("H132").Value = 17860246
$ws.Range("J132").Value = 4709.25
$ws.Range("L132").Value = 14127.75
$ws.Range("N132").Value = -19187.75
